$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 52
$ws.Range("F4").Value = 4621
$ws.Range("F5").Value = 1834
$ws.Range("F6").Value = 134
$ws.Range("F8").Value = 3101
$ws.Range("F11").Value = 257
$ws.Range("F12").Value = 622
$ws.Range("F13").Value = 531
$ws.Range("F14").Value = 523
$ws.Range("F15").Value = 373
$ws.Range("F17").Value = 1772
$ws.Range("F18").Value = 1322
$ws.Range("F19").Value = 122
$ws.Range("F20").Value = 1592
$ws.Range("F21").Value = 126
$ws.Range("F23").Value = 3
$ws.Range("F28").Value = 97
$ws.Range("F32").Value = 3772
$ws.Range("F33").Value = 757
$ws.Range("F35").Value = 638
$ws.Range("F37").Value = 1801

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 41

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 52
$ws.Range("F4").Value = 4621
$ws.Range("F5").Value = 1834
$ws.Range("F6").Value = 134
$ws.Range("F8").Value = 3101
$ws.Range("F11").Value = 257
$ws.Range("F12").Value = 622
$ws.Range("F13").Value = 531
$ws.Range("F14").Value = 523
$ws.Range("F16").Value = 373
$ws.Range("F18").Value = 1772
$ws.Range("F19").Value = 1322
$ws.Range("F20").Value = 122
$ws.Range("F21").Value = 1592
$ws.Range("F22").Value = 126
$ws.Range("F24").Value = 3
$ws.Range("F29").Value = 97
$ws.Range("F33").Value = 3772
$ws.Range("F34").Value = 41
$ws.Range("F35").Value = 757
$ws.Range("F37").Value = 638
$ws.Range("F39").Value = 1801
